$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto price/volume figures (cryptos.xlsx nightly
# GitHub Actions update). Most "Price" (column D) cells hold numeric-looking
# text (e.g. "1.00", "23.72", "0.0352") that must stay literal strings so
# trailing zeros / exact formatting survive - otherwise Excel's normal
# Range.Value auto-detection would coerce them into real numbers. For those
# cells we briefly force Text entry (NumberFormat "@") and then restore the
# cell's default style afterwards so no stray formatting is left behind.

$ws.Range("D2").Value = "62.409.79"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "3.010.92"
$ws.Range("E3").Value = "  -1.92%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "583.60"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "147.56"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.49%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").Value = "3.009.51"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -4.61%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.69"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.73%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.441"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000229"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "34.72"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -5.43%  "
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "3.507.93"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "62.403.54"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "3.010.56"
$ws.Range("E19").Value = "  -1.85%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "460.46"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.89%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.90"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.43%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.684"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.40%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.31"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.13%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.28"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -6.01%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "79.76"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -3.34%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.95"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("E30").Value = "  -1.08%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.36%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.09"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "26.97"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "0.0₃0785"
$ws.Range("E36").Value = "  -3.99%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.55%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.43%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "50.54"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "9.05"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.88"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -11.26%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "417.31"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("E43").Value = "  +0.90%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.274"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0352"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.773.26"
$ws.Range("E46").Value = "  -0.88%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "37.75"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -6.91%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "128.73"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -0.74%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "23.72"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -5.42%  "
